# Applies the "Updated symbol list" diff to the crypto price sheet.
# For each affected row this updates the Price (D) and Volume(1h) (E)
# columns with the freshly scraped values; rows 46 and 47 additionally
# swap their Coin (B) and Link (C) values (CoinbaseStockToken/BOLO order
# change) together with new Price/Volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map for every cell touched by the diff.
$edits = [ordered]@{
    "D2"  = "274.65";  "E2"  = "-1.61%"
    "D3"  = "26.64";   "E3"  = "-2.73%"
    "D4"  = "4.761";   "E4"  = "-0.91%"
    "D5"  = "0.06289"; "E5"  = "-1.00%"
    "D6"  = "6.917";   "E6"  = "-0.33%"
    "D7"  = "1.321";   "E7"  = "38.30%"
    "D8"  = "0.8684";  "E8"  = "-1.61%"
    "D9"  = "0.1591";  "E9"  = "7.31%"
    "D10" = "0.05036"; "E10" = "-3.30%"
    "D11" = "0.07484"; "E11" = "2.32%"
    "D12" = "0.02904"; "E12" = "-7.39%"
    "D13" = "0.09055"; "E13" = "-0.02%"
    "D14" = "0.001590"; "E14" = "1.57%"
    "D15" = "0.0006330"; "E15" = "1.14%"
    "D16" = "0.005925"; "E16" = "1.44%"
    "D17" = "3.454";   "E17" = "-0.26%"
    "D18" = "3.305";   "E18" = "-1.97%"
    "E19" = "-0.60%"
    "E20" = "0.86%"
    "D21" = "0.1318";  "E21" = "0.54%"
    "D22" = "3.922";   "E22" = "1.65%"
    "D23" = "0.04381"; "E23" = "1.49%"
    "D24" = "0.001168"; "E24" = "-1.02%"
    "D26" = "0.0001200"; "E26" = "0.03%"
    "D27" = "0.0001616"
    "D40" = "0.04066"; "E40" = "-0.55%"
    "D41" = "0.007098"; "E41" = "5.43%"
    "D42" = "0.1169";  "E42" = "0.51%"
    "D43" = "0.002020"; "E43" = "-11.38%"
    "D44" = "0.01123"; "E44" = "-10.12%"
    "D45" = "0.00005202"; "E45" = "-0.58%"
    "B46" = "CoinbaseStockToken"
    "C46" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
    "D46" = "0.02300"; "E46" = "2.15%"
    "B47" = "BOLO"
    "C47" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "D47" = "1.486";   "E47" = "-37.47%"
}

foreach ($addr in $edits.Keys) {
    $col = $addr.Substring(0, 1)
    $cell = $ws.Range($addr)

    if ($col -eq "D" -or $col -eq "E") {
        # Price and Volume values look numeric/percentage, but the sheet
        # stores them as plain text, so force a Text format first to
        # prevent Excel from auto-converting them to numbers.
        $cell.NumberFormat = "@"
    }

    $cell.Value = $edits[$addr]
}
